# Notebooks reran with CMap: "mean_score" columns became "mean_rank" columns,
# and per-sheet rows were re-sorted/re-scored accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: G3_effective
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("G3_effective")

$ws.Range("D1").Value = "mean_rank(G3)"
$ws.Range("F1").Value = "mean_rank(G4)"
$ws.Range("H1").Value = "mean_rank(SHH)"
$ws.Range("J1").Value = "mean_rank(SHH+p53)"

$ws.Range("A2").Value = 4
$ws.Range("D2").Value = 2.1
$ws.Range("F2").Value = 2.75
$ws.Range("H2").Value = 2
$ws.Range("J2").Value = 2

$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "bx-912"
$ws.Range("D3").Value = 10.4
$ws.Range("F3").Value = 8.75
$ws.Range("H3").Value = 116
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 223
$ws.Range("K3").Value = "['SBI-0645949.P001']"
$ws.Range("L3").Value = "SBI-0645949.P001"

$ws.Range("A4").Value = 44
$ws.Range("B4").Value = "abt-737"
$ws.Range("D4").Value = 32.05
$ws.Range("F4").Value = 32.75
$ws.Range("H4").Value = 33.41666666666666
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 33.83333333333334
$ws.Range("K4").Value = "[]"
$ws.Range("L4").Value = "NaN"

$ws.Range("A5").Value = 108
$ws.Range("B5").Value = "linsitinib"
$ws.Range("D5").Value = 84.34999999999999
$ws.Range("F5").Value = 88.75
$ws.Range("H5").Value = 87.25
$ws.Range("J5").Value = 86.5
$ws.Range("K5").Value = "['SBI-0646932.P001']"
$ws.Range("L5").Value = "SBI-0646932.P001"

# ---------------------------------------------------------------------------
# Sheet 2: G3_ineffective
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("G3_ineffective")

$ws.Range("D1").Value = "mean_rank(G3)"
$ws.Range("F1").Value = "mean_rank(G4)"
$ws.Range("H1").Value = "mean_rank(SHH)"
$ws.Range("J1").Value = "mean_rank(SHH+p53)"

$ws.Range("A2").Value = 135
$ws.Range("D2").Value = 156.95
$ws.Range("F2").Value = 167.625
$ws.Range("H2").Value = 163.5
$ws.Range("J2").Value = 184.6666666666667

$ws.Range("A3").Value = 149
$ws.Range("D3").Value = 228.05
$ws.Range("F3").Value = 217.25
$ws.Range("H3").Value = 214.6666666666667
$ws.Range("J3").Value = 189

# ---------------------------------------------------------------------------
# Sheet 3: notG3_effective
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("notG3_effective")

$ws.Range("D1").Value = "mean_rank(G3)"
$ws.Range("F1").Value = "mean_rank(G4)"
$ws.Range("H1").Value = "mean_rank(SHH)"
$ws.Range("J1").Value = "mean_rank(SHH+p53)"

$ws.Range("A2").Value = 12
$ws.Range("D2").Value = 7.9
$ws.Range("F2").Value = 8.5
$ws.Range("H2").Value = 7.666666666666667
$ws.Range("J2").Value = 7

$ws.Range("A3").Value = 64
$ws.Range("D3").Value = 34.3
$ws.Range("F3").Value = 48.25
$ws.Range("H3").Value = 55.25
$ws.Range("J3").Value = 42.5

# ---------------------------------------------------------------------------
# Sheet 4: SHH_effective
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SHH_effective")

$ws.Range("D1").Value = "mean_rank(G3)"
$ws.Range("F1").Value = "mean_rank(G4)"
$ws.Range("H1").Value = "mean_rank(SHH)"
$ws.Range("J1").Value = "mean_rank(SHH+p53)"

$ws.Range("A2").Value = 35
$ws.Range("B2").Value = "olaparib"
$ws.Range("D2").Value = 155.9
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 113.75
$ws.Range("H2").Value = 23
$ws.Range("J2").Value = 22.83333333333333

$ws.Range("A3").Value = 40
$ws.Range("B3").Value = "rucaparib"
$ws.Range("D3").Value = 113.4
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 21.75
$ws.Range("H3").Value = 15.75
$ws.Range("J3").Value = 16.16666666666667
